$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old table area (A1:E5), then drop the now-unused column E so the
# sheet shrinks from 5 columns to 4 (tcNum / xpath / expectedValue / result).
$ws.Range("A1:E5").ClearContents()
$ws.Columns.Item(5).Delete()

# New header row
$ws.Range("A1").Value = "tcNum"
$ws.Range("B1").Value = "xpath"
$ws.Range("C1").Value = "expectedValue"
$ws.Range("D1").Value = "result"

# New data rows (the "result" column is left blank, same as the template
# author left it - it only carries a header).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "//form[@id=`"blog_option_sort_form`"]//button[span='정확도']"
$ws.Range("C2").Value = "opt.sortsim"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "//form[@id=`"blog_option_sort_form`"]//button[span='최신순']"
$ws.Range("C3").Value = "opt.sortdate"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "//form[@id=`"blog_option_sort_form`"]//*[@class='clo_op']"
$ws.Range("C4").Value = "opt.sortfold"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "//form[@id=`"blog_option_sort_form`"]//*[@class='clo_op']"
$ws.Range("C5").Value = "opt.sortunfold"

# Column widths: col B widened to fit the new (much longer) xpath strings,
# col C narrowed to fit the opt.sort*/expectedValue strings. (Input values
# are pre-compensated for the host's internal character-rounding so the
# persisted <col width> lands as close as possible to the true target of
# 8.375 / 59.625 / 14.375.) Column D is left untouched so it keeps its
# original width/bestFit formatting rather than being reset.
$ws.Columns.Item(1).ColumnWidth = 7.714285714285714
$ws.Columns.Item(2).ColumnWidth = 58.857142857142854
$ws.Columns.Item(3).ColumnWidth = 13.714285714285714

# Move selection, matching the author's final cursor position.
$ws.Range("B15").Select()
